# Apply the latest cryptocurrency price/volume snapshot to the sheet.
#
# Column D ("Price") and column E ("Volume(1h)") values are stored as plain
# text in the workbook (e.g. "1.00", "539.30", "  +2.78%  "), not as numbers.
# To keep them as text (rather than have Excel silently reinterpret strings
# like "1.00" as the number 1) every new value is written with a leading
# apostrophe, which Excel treats as an explicit 'store as text' marker. The
# cell Style is then reset to "Normal" so that marker does not leave a visible
# quote-prefix / number-format change behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue 2 4 '60.404.94'   # Price: '60.222.67' -> '60.404.94'
Set-TextValue 2 5 '  +2.95%  '   # Volume(1h): '  +2.78%  ' -> '  +2.95%  '

# Row 3: Ethereum
Set-TextValue 3 4 '3.211.74'   # Price: '3.206.18' -> '3.211.74'
Set-TextValue 3 5 '  +1.50%  '   # Volume(1h): '  +1.38%  ' -> '  +1.50%  '

# Row 4: TetherUSD
Set-TextValue 4 4 '1.00'   # Price: '0.999' -> '1.00'
Set-TextValue 4 5 '  -0.05%  '   # Volume(1h): '  -0.14%  ' -> '  -0.05%  '

# Row 5: BNB
Set-TextValue 5 4 '539.32'   # Price: '539.30' -> '539.32'
Set-TextValue 5 5 '  +1.95%  '   # Volume(1h): '  +1.90%  ' -> '  +1.95%  '

# Row 6: Solana
Set-TextValue 6 4 '146.49'   # Price: '146.14' -> '146.49'
Set-TextValue 6 5 '  +4.74%  '   # Volume(1h): '  +4.77%  ' -> '  +4.74%  '

# Row 7: USDC
Set-TextValue 7 4 '0.999'   # Price: '1.00' -> '0.999'
Set-TextValue 7 5 '  -0.03%  '   # Volume(1h): '  -0.05%  ' -> '  -0.03%  '

# Row 8: XRP
Set-TextValue 8 5 '  -1.36%  '   # Volume(1h): '  -2.10%  ' -> '  -1.36%  '

# Row 9: Toncoin
Set-TextValue 9 4 '7.36'   # Price: '7.35' -> '7.36'
Set-TextValue 9 5 '  +0.96%  '   # Volume(1h): '  +0.92%  ' -> '  +0.96%  '

# Row 10: Dogecoin
Set-TextValue 10 5 '  +1.69%  '   # Volume(1h): '  +1.14%  ' -> '  +1.69%  '

# Row 11: Cardano
Set-TextValue 11 4 '0.434'   # Price: '0.433' -> '0.434'
Set-TextValue 11 5 '  -0.79%  '   # Volume(1h): '  -0.83%  ' -> '  -0.79%  '

# Row 12: WrappedliquidstakedEther2.0
Set-TextValue 12 4 '3.762.38'   # Price: '3.750.98' -> '3.762.38'
Set-TextValue 12 5 '  +1.47%  '   # Volume(1h): '  +1.12%  ' -> '  +1.47%  '

# Row 13: TRON
Set-TextValue 13 5 '  -1.94%  '   # Volume(1h): '  -1.88%  ' -> '  -1.94%  '

# Row 14: Avalanche
Set-TextValue 14 4 '25.99'   # Price: '25.88' -> '25.99'
Set-TextValue 14 5 '  +0.83%  '   # Volume(1h): '  +0.67%  ' -> '  +0.83%  '

# Row 15: ShibaInu
Set-TextValue 15 4 '0.0000173'   # Price: '0.0000172' -> '0.0000173'
Set-TextValue 15 5 '  +1.52%  '   # Volume(1h): '  +1.16%  ' -> '  +1.52%  '

# Row 16: WrappedBTC
Set-TextValue 16 4 '60.370.39'   # Price: '60.179.62' -> '60.370.39'
Set-TextValue 16 5 '  +2.85%  '   # Volume(1h): '  +2.60%  ' -> '  +2.85%  '

# Row 17: WrappedEther
Set-TextValue 17 4 '3.214.72'   # Price: '3.188.24' -> '3.214.72'
Set-TextValue 17 5 '  +1.67%  '   # Volume(1h): '  +1.51%  ' -> '  +1.67%  '

# Row 18: Polkadot
Set-TextValue 18 5 '  +0.52%  '   # Volume(1h): '  +0.61%  ' -> '  +0.52%  '

# Row 19: Chainlink
Set-TextValue 19 4 '13.30'   # Price: '13.28' -> '13.30'
Set-TextValue 19 5 '  +2.51%  '   # Volume(1h): '  +2.46%  ' -> '  +2.51%  '

# Row 20: Uniswap
Set-TextValue 20 4 '8.27'   # Price: '8.23' -> '8.27'
Set-TextValue 20 5 '  +2.09%  '   # Volume(1h): '  +1.66%  ' -> '  +2.09%  '

# Row 21: BitcoinCash
Set-TextValue 21 4 '372.63'   # Price: '371.56' -> '372.63'
Set-TextValue 21 5 '  -0.90%  '   # Volume(1h): '  -1.19%  ' -> '  -0.90%  '

# Row 22: Dai
Set-TextValue 22 4 '1.00'   # Price: '0.998' -> '1.00'
Set-TextValue 22 5 '  +0.02%  '   # Volume(1h): '  -0.19%  ' -> '  +0.02%  '

# Row 23: Polygon
Set-TextValue 23 5 '  -0.93%  '   # Volume(1h): '  -1.11%  ' -> '  -0.93%  '

# Row 24: Litecoin
Set-TextValue 24 4 '69.83'   # Price: '69.61' -> '69.83'
Set-TextValue 24 5 '  +0.06%  '   # Volume(1h): '  -0.21%  ' -> '  +0.06%  '

# Row 25: Kaspa
Set-TextValue 25 4 '0.171'   # Price: '0.170' -> '0.171'
Set-TextValue 25 5 '  +2.01%  '   # Volume(1h): '  +1.51%  ' -> '  +2.01%  '

# Row 26: InternetComputer(DFINITY)
Set-TextValue 26 5 '  +4.58%  '   # Volume(1h): '  +4.85%  ' -> '  +4.58%  '

# Row 27: Binance-PegBSC-USD
Set-TextValue 27 5 '  -0.21%  '   # Volume(1h): '  -0.06%  ' -> '  -0.21%  '

# Row 28: PEPE
Set-TextValue 28 4 '0.0₃0892'   # Price: '0.0₃0879' -> '0.0₃0892'
Set-TextValue 28 5 '  +3.31%  '   # Volume(1h): '  +1.91%  ' -> '  +3.31%  '

# Row 29: EthereumClassic
Set-TextValue 29 5 '  +0.40%  '   # Volume(1h): '  +0.32%  ' -> '  +0.40%  '

# Row 30: PancakeSwap
Set-TextValue 30 5 '  +0.89%  '   # Volume(1h): '  +0.81%  ' -> '  +0.89%  '

# Row 31: RenderToken
Set-TextValue 31 4 '6.15'   # Price: '6.12' -> '6.15'
Set-TextValue 31 5 '  +1.90%  '   # Volume(1h): '  +1.66%  ' -> '  +1.90%  '

# Row 32: NEARProtocol
Set-TextValue 32 5 '  +2.60%  '   # Volume(1h): '  +2.72%  ' -> '  +2.60%  '

# Row 33: Fetch.AI
Set-TextValue 33 4 '1.20'   # Price: '1.19' -> '1.20'
Set-TextValue 33 5 '  +4.25%  '   # Volume(1h): '  +3.60%  ' -> '  +4.25%  '

# Row 34: Aptos
Set-TextValue 34 4 '6.60'   # Price: '6.58' -> '6.60'
Set-TextValue 34 5 '  +4.53%  '   # Volume(1h): '  +4.36%  ' -> '  +4.53%  '

# Row 35: Monero
Set-TextValue 35 4 '158.65'   # Price: '158.83' -> '158.65'
Set-TextValue 35 5 '  +1.14%  '   # Volume(1h): '  +1.40%  ' -> '  +1.14%  '

# Row 36: ImmutableX
Set-TextValue 36 5 '  +3.07%  '   # Volume(1h): '  +3.52%  ' -> '  +3.07%  '

# Row 37: EnergySwap
Set-TextValue 37 4 '26.42'   # Price: '26.34' -> '26.42'
Set-TextValue 37 5 '  +5.67%  '   # Volume(1h): '  +4.98%  ' -> '  +5.67%  '

# Row 38: Maker
Set-TextValue 38 4 '2.799.37'   # Price: '2.797.10' -> '2.799.37'
Set-TextValue 38 5 '  +4.21%  '   # Volume(1h): '  +4.43%  ' -> '  +4.21%  '

# Row 39: VeChain
Set-TextValue 39 5 '  +8.92%  '   # Volume(1h): '  +9.28%  ' -> '  +8.92%  '

# Row 40: Hedera
Set-TextValue 40 4 '0.0712'   # Price: '0.0710' -> '0.0712'
Set-TextValue 40 5 '  +2.65%  '   # Volume(1h): '  +2.05%  ' -> '  +2.65%  '

# Row 41: Stacks
Set-TextValue 41 4 '1.70'   # Price: '1.69' -> '1.70'
Set-TextValue 41 5 '  +0.83%  '   # Volume(1h): '  +0.70%  ' -> '  +0.83%  '

# Row 42: Filecoin
Set-TextValue 42 4 '4.24'   # Price: '4.22' -> '4.24'
Set-TextValue 42 5 '  -0.96%  '   # Volume(1h): '  -1.38%  ' -> '  -0.96%  '

# Row 43: OKB
Set-TextValue 43 4 '39.97'   # Price: '39.98' -> '39.97'
Set-TextValue 43 5 '  +2.10%  '   # Volume(1h): '  +2.13%  ' -> '  +2.10%  '

# Row 44: Mantle
Set-TextValue 44 4 '0.721'   # Price: '0.720' -> '0.721'
Set-TextValue 44 5 '  -0.20%  '   # Volume(1h): '  -0.17%  ' -> '  -0.20%  '

# Row 45: Stellar
Set-TextValue 45 5 '  +2.30%  '   # Volume(1h): '  +1.63%  ' -> '  +2.30%  '

# Row 46: RenzoRestakedETH
Set-TextValue 46 4 '3.251.20'   # Price: '3.241.26' -> '3.251.20'
Set-TextValue 46 5 '  +1.42%  '   # Volume(1h): '  +1.10%  ' -> '  +1.42%  '

# Row 47: ONDO
Set-TextValue 47 4 '0.990'   # Price: '0.986' -> '0.990'
Set-TextValue 47 5 '  +1.20%  '   # Volume(1h): '  +0.83%  ' -> '  +1.20%  '

# Row 48: Cosmos
Set-TextValue 48 4 '6.17'   # Price: '6.16' -> '6.17'
Set-TextValue 48 5 '  -0.68%  '   # Volume(1h): '  -0.92%  ' -> '  -0.68%  '

# Row 49: InjectiveProtocol
Set-TextValue 49 4 '20.79'   # Price: '20.73' -> '20.79'
Set-TextValue 49 5 '  +3.72%  '   # Volume(1h): '  +3.54%  ' -> '  +3.72%  '

# Row 50: SuiNetwork
Set-TextValue 50 4 '0.801'   # Price: '0.798' -> '0.801'
Set-TextValue 50 5 '  +7.03%  '   # Volume(1h): '  +6.74%  ' -> '  +7.03%  '

# Row 51: FirstDigitalUSD
Set-TextValue 51 5 '  -0.04%  '   # Volume(1h): '  -0.13%  ' -> '  -0.04%  '
